# Apply the cell-content updates described by the target diff.
#
# Writing through Range.Value re-types each entry the way a user typing into
# Excel would: a string that parses as a plain number (e.g. "255.43", "1.00")
# is silently stored as a number, which would lose the original text
# formatting (trailing zeros, etc.) baked into these inline-string cells. To
# keep such values as text we prefix them with a literal leading apostrophe
# - Excel's standard "force text" entry marker - which is stripped from the
# stored value but keeps the cell type as text.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '41.243.45'
$ws.Range("E2").Value = '  +0.83%  '
$ws.Range("D3").Value = '2.188.92'
$ws.Range("E3").Value = '  -0.35%  '
$ws.Range("E4").Value = '  -0.04%  '
$ws.Range("D5").Value = "`'255.43"
$ws.Range("D6").Value = "`'0.628"
$ws.Range("E6").Value = '  +0.90%  '
$ws.Range("D7").Value = "`'68.32"
$ws.Range("E7").Value = '  -1.50%  '
$ws.Range("E8").Value = '  +0.03%  '
$ws.Range("D9").Value = "`'0.576"
$ws.Range("E9").Value = '  +6.03%  '
$ws.Range("D10").Value = "`'37.81"
$ws.Range("E10").Value = '  +3.77%  '
$ws.Range("D11").Value = "`'59.09"
$ws.Range("E11").Value = '  +2.86%  '
$ws.Range("D12").Value = "`'0.0937"
$ws.Range("E12").Value = '  -1.25%  '
$ws.Range("D13").Value = "`'7.14"
$ws.Range("E13").Value = '  +8.07%  '
$ws.Range("E14").Value = '  +0.78%  '
$ws.Range("D15").Value = '2.509.27'
$ws.Range("E15").Value = '  -0.47%  '
$ws.Range("D16").Value = "`'0.874"
$ws.Range("E16").Value = '  +5.11%  '
$ws.Range("D17").Value = "`'14.48"
$ws.Range("E17").Value = '  -1.18%  '
$ws.Range("D18").Value = '2.192.69'
$ws.Range("E18").Value = '  +0.08%  '
$ws.Range("D19").Value = '41.262.36'
$ws.Range("E19").Value = '  +1.04%  '
$ws.Range("D20").Value = '0.0₃0956'
$ws.Range("E20").Value = '  +1.51%  '
$ws.Range("D21").Value = "`'6.18"
$ws.Range("E21").Value = '  +2.18%  '
$ws.Range("D22").Value = "`'71.91"
$ws.Range("E22").Value = '  -0.76%  '
$ws.Range("D23").Value = "`'232.55"
$ws.Range("E23").Value = '  +1.22%  '
$ws.Range("E24").Value = '  +0.75%  '
$ws.Range("E25").Value = '  +10.21%  '
$ws.Range("D26").Value = "`'11.73"
$ws.Range("E26").Value = '  +21.39%  '
$ws.Range("D27").Value = "`'1.00"
$ws.Range("E27").Value = '  -0.02%  '
$ws.Range("E28").Value = '  +5.67%  '
$ws.Range("B29").Value = 'LEO'
$ws.Range("C29").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D29").Value = "`'3.75"
$ws.Range("E29").Value = '  -3.10%  '
$ws.Range("B30").Value = 'Toncoin'
$ws.Range("C30").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D30").Value = "`'2.17"
$ws.Range("E30").Value = '  -0.31%  '
$ws.Range("B31").Value = 'Monero'
$ws.Range("C31").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D31").Value = "`'168.80"
$ws.Range("E31").Value = '  -0.44%  '
$ws.Range("B32").Value = 'EthereumClassic'
$ws.Range("C32").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D32").Value = "`'20.67"
$ws.Range("E32").Value = '  +2.16%  '
$ws.Range("B33").Value = 'Kaspa'
$ws.Range("C33").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D33").Value = "`'0.117"
$ws.Range("E33").Value = '  -0.53%  '
$ws.Range("B34").Value = 'Hedera'
$ws.Range("C34").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D34").Value = "`'0.0755"
$ws.Range("E34").Value = '  +8.00%  '
$ws.Range("B35").Value = 'Stellar'
$ws.Range("C35").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D35").Value = "`'0.123"
$ws.Range("E35").Value = '  +0.03%  '
$ws.Range("B36").Value = 'InternetComputer(DFINITY)'
$ws.Range("C36").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D36").Value = "`'5.47"
$ws.Range("E36").Value = '  +6.70%  '
$ws.Range("B37").Value = 'InjectiveProtocol'
$ws.Range("C37").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D37").Value = "`'26.32"
$ws.Range("E37").Value = '  +9.52%  '
$ws.Range("B38").Value = 'Filecoin'
$ws.Range("C38").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D38").Value = "`'4.63"
$ws.Range("E38").Value = '  +1.16%  '
$ws.Range("B39").Value = 'RenderToken'
$ws.Range("C39").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D39").Value = "`'4.14"
$ws.Range("E39").Value = '  +8.36%  '
$ws.Range("B40").Value = 'VeChain'
$ws.Range("C40").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D40").Value = "`'0.0300"
$ws.Range("E40").Value = '  +10.84%  '
$ws.Range("D41").Value = "`'12.62"
$ws.Range("E41").Value = '  +21.14%  '
$ws.Range("B42").Value = 'LidoDAOToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D42").Value = "`'2.20"
$ws.Range("E42").Value = '  -2.70%  '
$ws.Range("B43").Value = 'THORChain'
$ws.Range("C43").Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range("D43").Value = "`'5.67"
$ws.Range("E43").Value = '  -1.81%  '
$ws.Range("B44").Value = 'FTXToken'
$ws.Range("C44").Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range("D44").Value = "`'5.15"
$ws.Range("E44").Value = '  +5.81%  '
$ws.Range("B45").Value = 'MultiversX'
$ws.Range("C45").Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range("D45").Value = "`'64.33"
$ws.Range("E45").Value = '  +3.16%  '
$ws.Range("B46").Value = 'Algorand'
$ws.Range("C46").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D46").Value = "`'0.201"
$ws.Range("E46").Value = '  +4.95%  '
$ws.Range("B47").Value = 'FraxShare'
$ws.Range("C47").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D47").Value = "`'8.64"
$ws.Range("E47").Value = '  +0.79%  '
$ws.Range("B48").Value = 'Cronos'
$ws.Range("C48").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D48").Value = "`'0.101"
$ws.Range("E48").Value = '  +3.09%  '
$ws.Range("B49").Value = 'BinanceUSD'
$ws.Range("C49").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range("D49").Value = "`'1.00"
$ws.Range("E49").Value = '  +0.28%  '
$ws.Range("B50").Value = 'ARBITRUM'
$ws.Range("C50").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D50").Value = "`'1.14"
$ws.Range("E50").Value = '  +4.99%  '
$ws.Range("B51").Value = 'TrustWalletToken'
$ws.Range("C51").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D51").Value = "`'1.17"
$ws.Range("E51").Value = '  +1.15%  '
